$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.462.01"
$ws.Range("E2").Value = "  -1.30%  "
$ws.Range("D3").Value = "2.338.18"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "314.30"
$ws.Range("E5").Value = "  -4.63%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "107.26"
$ws.Range("E6").Value = "  +1.29%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.630"
$ws.Range("E7").Value = "  -1.78%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.614"
$ws.Range("E9").Value = "  -6.94%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.77"
$ws.Range("E10").Value = "  -2.03%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0922"
$ws.Range("E11").Value = "  -1.74%  "
$ws.Range("E12").Value = "  -2.28%  "
$ws.Range("E13").Value = "  +0.01%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.987"
$ws.Range("E14").Value = "  -5.57%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.86"
$ws.Range("E15").Value = "  -7.49%  "
$ws.Range("D16").Value = "2.694.49"
$ws.Range("E16").Value = "  -1.81%  "
$ws.Range("D17").Value = "2.354.74"
$ws.Range("E17").Value = "  -0.94%  "
$ws.Range("D18").Value = "42.461.95"
$ws.Range("E18").Value = "  -1.57%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.56"
$ws.Range("E19").Value = "  -5.53%  "
$ws.Range("E20").Value = "  -2.91%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "75.77"
$ws.Range("E21").Value = "  -1.34%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.63"
$ws.Range("E22").Value = "  -1.80%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "255.54"
$ws.Range("E23").Value = "  -8.35%  "
$ws.Range("E24").Value = "  -5.50%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.32"
$ws.Range("E25").Value = "  -3.01%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  -0.05%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.32"
$ws.Range("E27").Value = "  -3.86%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "22.75"
$ws.Range("E28").Value = "  -2.01%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "173.62"
$ws.Range("E30").Value = "  -0.71%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "36.45"
$ws.Range("E31").Value = "  -4.23%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0885"
$ws.Range("E32").Value = "  -4.61%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.04"
$ws.Range("E33").Value = "  +2.80%  "
$ws.Range("E34").Value = "  -9.92%  "
$ws.Range("E35").Value = "  +17.41%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.131"
$ws.Range("E36").Value = "  -2.06%  "
$ws.Range("E37").Value = "  -5.88%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0360"
$ws.Range("E38").Value = "  -2.20%  "
$ws.Range("E39").Value = "  -10.62%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.63"
$ws.Range("E40").Value = "  -6.98%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.235"
$ws.Range("E41").Value = "  -0.28%  "
$ws.Range("B42").Value = "MultiversX"
$ws.Range("C42").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "70.37"
$ws.Range("E42").Value = "  +0.51%  "
$ws.Range("B43").Value = "ARBITRUM"
$ws.Range("C43").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.46"
$ws.Range("E43").Value = "  -7.96%  "
$ws.Range("E44").Value = "  -0.17%  "
$ws.Range("E45").Value = "  -4.86%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "110.90"
$ws.Range("E46").Value = "  -10.37%  "
$ws.Range("B47").Value = "FraxShare"
$ws.Range("C47").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.14"
$ws.Range("E47").Value = "  -3.05%  "
$ws.Range("B48").Value = "BitcoinSV"
$ws.Range("C48").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "84.98"
$ws.Range("E48").Value = "  -10.11%  "
$ws.Range("B49").Value = "THORChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.42"
$ws.Range("E49").Value = "  -2.15%  "
$ws.Range("B50").Value = "ordi"
$ws.Range("C50").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "74.01"
$ws.Range("E50").Value = "  +1.64%  "
$ws.Range("B51").Value = "TrustWalletToken"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.27"
$ws.Range("E51").Value = "  -3.30%  "
